# Trade #6 closed at 2026-02-17 15:13:55 - unknown UNKNOWN +0.000%

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Summary sheet: refresh aggregate stats now that trade #6 has closed
# ---------------------------------------------------------------------
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B3").Value = 1199.89   # Current Capital
$summary.Range("B4").Value = -0.11     # Total P&L $
$summary.Range("B5").Value = -0.37     # Total P&L %
$summary.Range("B6").Value = 6         # Total Trades
$summary.Range("B8").Value = 3         # Losing Trades
$summary.Range("B9").Value = 33.33     # Win Rate %

# ---------------------------------------------------------------------
# Strategy Status sheet: MarketMaking row (row 4) stats
# ---------------------------------------------------------------------
$status = $wb.Worksheets.Item("Strategy Status")
$status.Range("C4").Value = 99.89      # Capital
$status.Range("D4").Value = 6          # Trades
$status.Range("E4").Value = -0.11      # P&L $
$status.Range("F4").Value = -0.11      # P&L %
$status.Range("G4").Value = 33.33      # Win Rate %

# ---------------------------------------------------------------------
# All Trades + MarketMaking sheets: append trade #6 as new row 7
# ---------------------------------------------------------------------
$sheetNames = @("All Trades", "MarketMaking")
foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $row = 7
    $ws.Cells.Item($row, 1).Value = 6

    # Date/time columns are stored as plain text in this workbook (not
    # Excel date serials) - force Text format first so "2026-02-17"
    # isn't auto-converted into a date value.
    $ws.Cells.Item($row, 2).NumberFormat = "@"
    $ws.Cells.Item($row, 2).Value = "2026-02-17"
    $ws.Cells.Item($row, 3).Value = "15:13:49"

    $ws.Cells.Item($row, 4).Value = "MarketMaking"
    $ws.Cells.Item($row, 5).Value = "DOWN"
    $ws.Cells.Item($row, 6).Value = 0.158562
    $ws.Cells.Item($row, 7).Value = 0.06
    $ws.Cells.Item($row, 8).Value = "CLOSED"
    $ws.Cells.Item($row, 9).Value = -62.16
    $ws.Cells.Item($row, 10).Value = -0.1
    $ws.Cells.Item($row, 11).Value = 99.89
    $ws.Cells.Item($row, 12).Value = 0
    $ws.Cells.Item($row, 13).Value = 0
    $ws.Cells.Item($row, 14).Value = 0.6
    $ws.Cells.Item($row, 15).Value = "Normal spread capture: 19600 bps"
    $ws.Cells.Item($row, 16).Value = "early_exit"
    $ws.Cells.Item($row, 17).Value = 0.15
}
